# Add a new "08-ago" column (AU) after the existing "07-ago" column (AT),
# populating the header and the daily values for each of the 17 data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("AU1").Value = "08-ago"

# Data values (row -> AU value)
$ws.Range("AU2").Value = 0
$ws.Range("AU3").Value = 13.00976649845494
$ws.Range("AU4").Value = 18.670726704322401
$ws.Range("AU5").Value = 15.543866453957417
$ws.Range("AU6").Value = 0
$ws.Range("AU7").Value = 8.4651412514770641
$ws.Range("AU8").Value = 14.613141895780513
$ws.Range("AU9").Value = 11.426605606748559
$ws.Range("AU10").Value = 17.512875661830336
$ws.Range("AU11").Value = 14.079554146773756
$ws.Range("AU12").Value = 0
$ws.Range("AU13").Value = 6.7693726043457261
$ws.Range("AU14").Value = 0
$ws.Range("AU15").Value = 0
$ws.Range("AU16").Value = 11.335939523265097
$ws.Range("AU17").Value = 0
$ws.Range("AU18").Value = 0

# The active selection shifted one column to the right (AV5 -> AW5) as a
# side effect of the new column being added.
$ws.Range("AW5").Select()
